# Apply the scraped-data refresh described in the commit:
# "Update gh-pages to output generated at 456a3b4"
#
# This updates "want to go" counts (F) / min-price (G) figures that
# naturally drifted since the last scrape, plus a handful of rows whose
# underlying events were replaced/renamed on bilibili's event listing.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "展览" (Exhibitions) - F-column "want to go" count refresh
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 323
$ws.Range("F4").Value  = 2999
$ws.Range("F5").Value  = 79
$ws.Range("F7").Value  = 2331
$ws.Range("F8").Value  = 1705
$ws.Range("F9").Value  = 61
$ws.Range("F10").Value = 862
$ws.Range("F13").Value = 1
$ws.Range("F17").Value = 7134
$ws.Range("F19").Value = 7279
$ws.Range("F22").Value = 5554
$ws.Range("F27").Value = 194
$ws.Range("F28").Value = 1918
$ws.Range("F30").Value = 309
$ws.Range("F33").Value = 292
$ws.Range("F35").Value = 2443
$ws.Range("F36").Value = 1233
$ws.Range("F37").Value = 2789
$ws.Range("F38").Value = 39
$ws.Range("F41").Value = 398

# ---------------------------------------------------------------
# Sheet "演出" (Performances)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value  = 220
$ws.Range("F12").Value = 319
# This show is now sold out - min price column switches from a
# numeric price to the literal text "已售罄" (sold out).
$ws.Range("G12").Value = "已售罄"
$ws.Range("F18").Value = 67
$ws.Range("F20").Value = 7

# ---------------------------------------------------------------
# Sheet "本地生活" (Local life)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 79

# ---------------------------------------------------------------
# Sheet "全部类型" (All types) - union of the other sheets, plus a
# few rows whose underlying events were swapped out entirely.
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 323
$ws.Range("F5").Value  = 2999
$ws.Range("F6").Value  = 2331
$ws.Range("F7").Value  = 1705
$ws.Range("F8").Value  = 61
$ws.Range("F9").Value  = 862
$ws.Range("F13").Value = 79
$ws.Range("F16").Value = 220
$ws.Range("F19").Value = 7134

# Row 20: "...孟宇专场活动" event replaced by "...一周年盛典"
# (which previously sat in row 21).
$ws.Range("C20").Value = "北京·第十三届GOJO超次元动漫游戏嘉年华·一周年盛典"
$ws.Range("E20").Value = "2024.06.01 10:00-06.02 17:00"
$ws.Range("F20").Value = 7279
$ws.Range("G20").Value = 6.6
$ws.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=83827"
$ws.Range("I20").Value = "//i2.hdslb.com/bfs/openplatform/202405/T7pPJ1NM1715065435932.jpeg"

# Row 21: replaced by a brand new event "北京·银魂only2.0".
$ws.Range("C21").Value = "北京·银魂only2.0"
$ws.Range("D21").Value = "太平庄中街西端 北京天通苑黄河京都会议中心"
$ws.Range("E21").Value = "2024.06.01 10:00-06.01 17:00"
$ws.Range("F21").Value = 5
$ws.Range("G21").Value = 129
$ws.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=85154"
$ws.Range("I21").Value = "//i2.hdslb.com/bfs/openplatform/202405/7hbieM0H1714994580295.jpeg"

$ws.Range("F23").Value = 5554

# Row 25: "紫罗兰永恒花园交响音乐会" event replaced by "亦创·梦次元动漫游戏展1st"
# (which previously sat in row 26).
$ws.Range("C25").Value = "北京·亦创·梦次元动漫游戏展1st"
$ws.Range("D25").Value = "亦庄荣昌东街6号 北京亦创国际会展中心"
$ws.Range("E25").Value = "2024.06.08 09:30-06.08 17:00"
$ws.Range("F25").Value = 3495
$ws.Range("G25").Value = 80
$ws.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=84015"
$ws.Range("I25").Value = "//i1.hdslb.com/bfs/openplatform/202404/UfpmzLsm1712649924888.jpeg"

# Row 26: replaced by a brand new event "北京·嘉品次元派对·免费展会".
# Leading apostrophe forces this date-shaped string to stay plain text
# (matching the rest of column B) instead of being parsed as a date;
# resetting the style afterwards drops the quote-prefix formatting so
# the cell keeps the sheet's default (unstyled) look, like its neighbors.
$ws.Range("B26").Value = "'2024-06-09"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "北京·嘉品次元派对·免费展会"
$ws.Range("D26").Value = "东坝中路38号 北京金隅嘉品Mall中庭"
$ws.Range("E26").Value = "2024.06.09 14:00-06.10 20:30"
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=85726"
$ws.Range("I26").Value = "//i2.hdslb.com/bfs/openplatform/202405/v0azRLWZ1715829225052.jpeg"

$ws.Range("F30").Value = 1918
$ws.Range("F33").Value = 309
$ws.Range("F36").Value = 292
$ws.Range("F38").Value = 2443
$ws.Range("F39").Value = 1233
$ws.Range("F40").Value = 67
$ws.Range("F41").Value = 2789
$ws.Range("F42").Value = 39
$ws.Range("F45").Value = 398
